$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Price (D) values are textual (e.g. "30.793.19" / "0.9967") in the source data,
# so force the cell format to Text before assigning -- otherwise Excel
# auto-coerces plain-decimal-looking strings into floating point numbers.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '30.793.19'
$ws.Range("E2").Value = '  +0.85%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.930.08'
$ws.Range("E3").Value = '  +2.37%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.9967'
$ws.Range("E4").Value = '  -0.49%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '250.03'
$ws.Range("E5").Value = '  +2.28%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.5946'
$ws.Range("E6").Value = '  +26.28%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.9962'
$ws.Range("E7").Value = '  -0.53%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3117'
$ws.Range("E8").Value = '  +7.38%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '24.20'
$ws.Range("E9").Value = '  +8.07%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.06693'
$ws.Range("E10").Value = '  +3.29%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.7894'
$ws.Range("E11").Value = '  +7.21%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '101.38'
$ws.Range("E12").Value = '  +6.26%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.07923'
$ws.Range("E13").Value = '  +2.20%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '1.908.85'
$ws.Range("E14").Value = '  +1.24%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '5.338'
$ws.Range("E15").Value = '  +3.13%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '284.34'
$ws.Range("E16").Value = '  +0.53%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '30.694.14'
$ws.Range("E17").Value = '  +0.54%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '13.67'
$ws.Range("E18").Value = '  +4.88%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.000007662'
$ws.Range("E19").Value = '  +2.48%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '5.502'
$ws.Range("E20").Value = '  +4.72%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '0.9976'
$ws.Range("E21").Value = '  -0.36%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '2.152.89'
$ws.Range("E22").Value = '  +1.23%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '0.9970'
$ws.Range("E23").Value = '  -0.48%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '6.591'
$ws.Range("E24").Value = '  +5.45%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '9.332'
$ws.Range("E25").Value = '  +2.15%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '165.06'
$ws.Range("E26").Value = '  +0.72%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '19.52'
$ws.Range("E27").Value = '  +3.76%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '1.975'
$ws.Range("E28").Value = '  +4.31%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '0.1073'
$ws.Range("E29").Value = '  +10.45%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '1.352'
$ws.Range("E30").Value = '  +1.12%  '
$ws.Range("E31").Value = '  +3.36%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '4.468'
$ws.Range("E32").Value = '  +4.48%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '4.298'
$ws.Range("E33").Value = '  +4.30%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.04961'
$ws.Range("E34").Value = '  +1.94%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.157'
$ws.Range("E35").Value = '  +2.74%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.7167'
$ws.Range("E36").Value = '  +3.46%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '2.769'
$ws.Range("E37").Value = '  +2.33%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.01948'
$ws.Range("E38").Value = '  +2.50%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '2.918'
$ws.Range("E39").Value = '  +3.13%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.4423'
$ws.Range("E42").Value = '  +4.05%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '2.005'
$ws.Range("E43").Value = '  +0.21%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.8415'
$ws.Range("E44").Value = '  +2.27%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.9983'
$ws.Range("E45").Value = '  -0.26%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '10.06'
$ws.Range("E46").Value = '  +5.14%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '102.06'
$ws.Range("E47").Value = '  +0.87%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '7.204'
$ws.Range("E48").Value = '  +3.47%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '35.67'
$ws.Range("E49").Value = '  +1.15%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.4092'
$ws.Range("E50").Value = '  +3.99%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '910.28'
$ws.Range("E51").Value = '  -0.18%  '

# Rows 40/41: ranking reorder swaps FraxShare <-> Aave (with updated price/volume)
$ws.Range("B40").Value = 'FraxShare'
$ws.Range("C40").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '6.431'
$ws.Range("E40").Value = '  +4.26%  '
$ws.Range("B41").Value = 'Aave'
$ws.Range("C41").Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '77.21'
$ws.Range("E41").Value = '  +2.81%  '
